# Apply weekly update: prepend 3 new price rows (week of 2021-10-07, serial 44476)
# for "Murcott" quality grades Especial/Primera/Segunda, pushing existing data down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above row 226 (existing rows 226.. shift down to 229..)
$ws.Rows.Item(226).Insert()
$ws.Rows.Item(226).Insert()
$ws.Rows.Item(226).Insert()

# Common / static column values shared by all rows in this block
$mercadoId   = 2
$mercado     = "Comercializadora del Agro de Limarí"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102004
$categoria   = "Mandarina"
$variedad    = "Murcott"
$unidad      = "$/bandeja 10 kilos"
$origen      = "Provincia de Limarí"
$fecha       = 44476

function Set-PrecioRow($Row, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg, $KgUnidad) {
    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $tipo
    $ws.Cells.Item($Row, 7).Value  = $productoId
    $ws.Cells.Item($Row, 8).Value  = $producto
    $ws.Cells.Item($Row, 9).Value  = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

Set-PrecioRow 226 "Especial" 500 4500 5000 4750 475 10
Set-PrecioRow 227 "Primera"  700 3500 4000 3750 375 10
Set-PrecioRow 228 "Segunda"  500 2500 3000 2750 275 10
